$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.772.58"
$ws.Range("E2").Value = "  +4.62%  "

$ws.Range("D3").Value = "1.612.23"
$ws.Range("E3").Value = "  +3.97%  "

$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'213.94"
$ws.Range("E5").Value = "  +1.51%  "

$ws.Range("D6").Value = "'0.516"
$ws.Range("E6").Value = "  +6.96%  "

$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("D8").Value = "'26.72"
$ws.Range("E8").Value = "  +11.72%  "

$ws.Range("D9").Value = "'0.250"
$ws.Range("E9").Value = "  +3.38%  "

$ws.Range("D10").Value = "'0.0599"
$ws.Range("E10").Value = "  +2.70%  "

$ws.Range("D11").Value = "'0.0911"
$ws.Range("E11").Value = "  +2.54%  "

$ws.Range("D12").Value = "1.844.76"
$ws.Range("E12").Value = "  +4.08%  "

$ws.Range("D13").Value = "1.620.17"
$ws.Range("E13").Value = "  +4.37%  "

$ws.Range("D14").Value = "29.810.02"
$ws.Range("E14").Value = "  +4.82%  "

$ws.Range("D15").Value = "'0.539"
$ws.Range("E15").Value = "  +5.87%  "

$ws.Range("E16").Value = "  +3.81%  "

$ws.Range("D17").Value = "'244.83"
$ws.Range("E17").Value = "  +7.19%  "

$ws.Range("D18").Value = "'63.54"
$ws.Range("E18").Value = "  +4.18%  "

$ws.Range("E19").Value = "  +4.38%  "

$ws.Range("D20").Value = "0.0₃0695"
$ws.Range("E20").Value = "  +3.42%  "

$ws.Range("D21").Value = "'0.996"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").Value = "'4.05"
$ws.Range("E22").Value = "  +4.11%  "

$ws.Range("D23").Value = "'9.28"
$ws.Range("E23").Value = "  +4.26%  "

$ws.Range("D24").Value = "'2.10"
$ws.Range("E24").Value = "  +4.19%  "

$ws.Range("D25").Value = "'155.72"
$ws.Range("E25").Value = "  +3.27%  "

$ws.Range("D26").Value = "'15.37"
$ws.Range("E26").Value = "  +4.32%  "

$ws.Range("E27").Value = "  +5.52%  "

$ws.Range("D28").Value = "'6.41"
$ws.Range("E28").Value = "  +3.23%  "

$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("E30").Value = "  +1.39%  "

$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("E32").Value = "  +3.09%  "

$ws.Range("D33").Value = "1.439.75"
$ws.Range("E33").Value = "  +4.15%  "

$ws.Range("D34").Value = "'3.11"
$ws.Range("E34").Value = "  +3.72%  "

$ws.Range("E35").Value = "  -1.27%  "

$ws.Range("D36").Value = "'2.84"
$ws.Range("E36").Value = "  +10.77%  "

$ws.Range("E37").Value = "  +2.95%  "

$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("E39").Value = "  +3.20%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.538"
$ws.Range("E40").Value = "  +5.29%  "

$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").Value = "'55.53"
$ws.Range("E41").Value = "  +28.76%  "

$ws.Range("D42").Value = "'1.95"
$ws.Range("E42").Value = "  +1.56%  "

$ws.Range("D43").Value = "'0.798"
$ws.Range("E43").Value = "  +3.47%  "

$ws.Range("D44").Value = "'0.996"
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("D45").Value = "'0.0469"
$ws.Range("E45").Value = "  +2.43%  "

$ws.Range("D46").Value = "'66.18"
$ws.Range("E46").Value = "  +7.19%  "

$ws.Range("D47").Value = "'5.34"
$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("D48").Value = "1.754.28"
$ws.Range("E48").Value = "  +4.16%  "

$ws.Range("D49").Value = "'86.97"
$ws.Range("E49").Value = "  +2.15%  "

$ws.Range("E50").Value = "  -4.44%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0520"
$ws.Range("E51").Value = "  +1.86%  "
